$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("cofi")
$ws.Activate()

# --- Prepare "New" section header cell format at A15 (copy fill-only style from A4, strip border) ---
# (value is assigned later so the shared-string table keeps the same ordering as the source edit)
$ws.Range("A4").Copy()
$ws.Range("A15").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A15").Borders.LineStyle = -4142 # xlLineStyleNone

# --- Row 17: base / example hierarchy row ---
$ws.Range("A17").Value = 0
$ws.Range("B17").Value = "Customer Care"
$ws.Range("C17").Value = 0
$ws.Range("D17").Value = "cases.csv"
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = "CaseNumber"
$ws.Range("G17").Value = "C0000000000000000000000000000000"
$ws.Range("H17").Value = "OI000000000000000000000000000000"
$ws.Range("I17").Value = "FI000000000000000000000000000000"

# --- Row 18 ---
$ws.Range("A18").Value = 0
$ws.Range("B18").Value = "Customer Care"
$ws.Range("C18").Value = 0
$ws.Range("D18").Value = "cases.csv"
$ws.Range("E18").Value = 1
$ws.Range("F18").Value = "CreatedDate"
$ws.Range("G18").Value = "C0000000000000000000000000000000"
$ws.Range("H18").Value = "OI000000000000000000000000000000"
$ws.Range("I18").Value = "FI000000000000000000000000000001"

# --- Row 19 ---
$ws.Range("A19").Value = 0
$ws.Range("B19").Value = "Customer Care"
$ws.Range("C19").Value = 0
$ws.Range("D19").Value = "cases.csv"
$ws.Range("E19").Value = 2
$ws.Range("F19").Value = "Description"
$ws.Range("G19").Value = "C0000000000000000000000000000000"
$ws.Range("H19").Value = "OI000000000000000000000000000000"
$ws.Range("I19").Value = "FI000000000000000000000000000002"

# --- Row 20 ---
$ws.Range("A20").Value = 0
$ws.Range("B20").Value = "Customer Care"
$ws.Range("C20").Value = 0
$ws.Range("D20").Value = "cases.csv"
$ws.Range("E20").Value = 3
$ws.Range("F20").Value = "Id"
$ws.Range("G20").Value = "C0000000000000000000000000000000"
$ws.Range("H20").Value = "OI000000000000000000000000000000"
$ws.Range("I20").Value = "FI000000000000000000000000000003"

# --- Row 22: column header labels for the UUID / COFI hierarchy table ---
$ws.Range("A22").Value = "L1 Property"
$ws.Range("B22").Value = "L1 Display Name"
$ws.Range("C22").Value = "L2 Property"
$ws.Range("D22").Value = "L2 Display Name"
$ws.Range("E22").Value = "L3 Property"
$ws.Range("F22").Value = "L3 Display Name"
$ws.Range("G22").Value = "L1 Code"
$ws.Range("H22").Value = "L2 Code"
$ws.Range("I22").Value = "L3 Code"

# --- Now set the "New" section header text (registered last in shared strings) ---
$ws.Range("A15").Value = "New"

# --- Column widths: columns A (1) through L (12) set to stored width 17.5 ---
$ws.Range("A1:L1").ColumnWidth = 16.666666666666668

# --- Update the active cell selection to match the author's final cursor position ---
$null = $ws.Range("B18").Select()
